# Update the NATMI ligand-receptor pair statistics (Gnas-Lhcgr) with the
# recomputed TPM-based values. Columns: G=Ligand avg expr, H=Ligand total
# expr, I/J=Ligand specificity (avg/total), K=Receptor-expressing cells,
# L=Receptor detection rate, M=Receptor avg expr, N=Receptor total expr,
# O/P=Receptor specificity (avg/total), Q/R=Edge expr weight (avg/total),
# S/T=Edge specificity (avg/total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 189.0573523333333
$ws.Cells.Item(2, 8).Value = 567.172057
$ws.Cells.Item(2, 9).Value = 0.1182556374491171
$ws.Cells.Item(2, 10).Value = 0.1182556374491171
$ws.Cells.Item(2, 14).Value = 0.7920560000000001
$ws.Cells.Item(2, 15).Value = 0.4210262213814106
$ws.Cells.Item(2, 16).Value = 0.4210262213814106
$ws.Cells.Item(2, 17).Value = 49.91467008657689
$ws.Cells.Item(2, 18).Value = 449.232030779192
$ws.Cells.Item(2, 19).Value = 0.0497887241922518
$ws.Cells.Item(2, 20).Value = 0.0497887241922518

$ws.Cells.Item(3, 7).Value = 189.0573523333333
$ws.Cells.Item(3, 8).Value = 567.172057
$ws.Cells.Item(3, 9).Value = 0.1182556374491171
$ws.Cells.Item(3, 10).Value = 0.1182556374491171
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.363065
$ws.Cells.Item(3, 14).Value = 1.089195
$ws.Cells.Item(3, 15).Value = 0.5789737786185895
$ws.Cells.Item(3, 16).Value = 0.5789737786185893
$ws.Cells.Item(3, 17).Value = 68.64010762490166
$ws.Cells.Item(3, 18).Value = 617.760968624115
$ws.Cells.Item(3, 19).Value = 0.0684669132568653
$ws.Cells.Item(3, 20).Value = 0.06846691325686528

$ws.Cells.Item(4, 7).Value = 930.1503093333332
$ws.Cells.Item(4, 9).Value = 0.5818103152093762
$ws.Cells.Item(4, 10).Value = 0.5818103152093762
$ws.Cells.Item(4, 14).Value = 0.7920560000000001
$ws.Cells.Item(4, 15).Value = 0.4210262213814106
$ws.Cells.Item(4, 16).Value = 0.4210262213814106
$ws.Cells.Item(4, 17).Value = 245.5770444697742
$ws.Cells.Item(4, 19).Value = 0.2449573985733311
$ws.Cells.Item(4, 20).Value = 0.2449573985733311

$ws.Cells.Item(5, 7).Value = 930.1503093333332
$ws.Cells.Item(5, 9).Value = 0.5818103152093762
$ws.Cells.Item(5, 10).Value = 0.5818103152093762
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.363065
$ws.Cells.Item(5, 14).Value = 1.089195
$ws.Cells.Item(5, 15).Value = 0.5789737786185895
$ws.Cells.Item(5, 16).Value = 0.5789737786185893
$ws.Cells.Item(5, 17).Value = 337.7050220581066
$ws.Cells.Item(5, 18).Value = 3039.345198522959
$ws.Cells.Item(5, 19).Value = 0.3368529166360451
$ws.Cells.Item(5, 20).Value = 0.3368529166360451

$ws.Cells.Item(6, 7).Value = 420.6651306666666
$ws.Cells.Item(6, 8).Value = 1261.995392
$ws.Cells.Item(6, 9).Value = 0.2631266256807295
$ws.Cells.Item(6, 10).Value = 0.2631266256807295
$ws.Cells.Item(6, 14).Value = 0.7920560000000001
$ws.Cells.Item(6, 15).Value = 0.4210262213814106
$ws.Cells.Item(6, 16).Value = 0.4210262213814106
$ws.Cells.Item(6, 17).Value = 111.0634469117724
$ws.Cells.Item(6, 18).Value = 999.571022205952
$ws.Cells.Item(6, 19).Value = 0.1107832089551984
$ws.Cells.Item(6, 20).Value = 0.1107832089551984

$ws.Cells.Item(7, 7).Value = 420.6651306666666
$ws.Cells.Item(7, 8).Value = 1261.995392
$ws.Cells.Item(7, 9).Value = 0.2631266256807295
$ws.Cells.Item(7, 10).Value = 0.2631266256807295
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.363065
$ws.Cells.Item(7, 14).Value = 1.089195
$ws.Cells.Item(7, 15).Value = 0.5789737786185895
$ws.Cells.Item(7, 16).Value = 0.5789737786185893
$ws.Cells.Item(7, 17).Value = 152.7287856654933
$ws.Cells.Item(7, 18).Value = 1374.55907098944
$ws.Cells.Item(7, 19).Value = 0.1523434167255311
$ws.Cells.Item(7, 20).Value = 0.1523434167255311

$ws.Cells.Item(8, 7).Value = 58.84466766666667
$ws.Cells.Item(8, 8).Value = 176.534003
$ws.Cells.Item(8, 9).Value = 0.03680742166077718
$ws.Cells.Item(8, 10).Value = 0.03680742166077718
$ws.Cells.Item(8, 14).Value = 0.7920560000000001
$ws.Cells.Item(8, 15).Value = 0.4210262213814106
$ws.Cells.Item(8, 16).Value = 0.4210262213814106
$ws.Cells.Item(8, 17).Value = 15.53609069779645
$ws.Cells.Item(8, 18).Value = 139.824816280168
$ws.Cells.Item(8, 19).Value = 0.0154968896606293
$ws.Cells.Item(8, 20).Value = 0.0154968896606293

$ws.Cells.Item(9, 7).Value = 58.84466766666667
$ws.Cells.Item(9, 8).Value = 176.534003
$ws.Cells.Item(9, 9).Value = 0.03680742166077718
$ws.Cells.Item(9, 10).Value = 0.03680742166077718
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.363065
$ws.Cells.Item(9, 14).Value = 1.089195
$ws.Cells.Item(9, 15).Value = 0.5789737786185895
$ws.Cells.Item(9, 16).Value = 0.5789737786185893
$ws.Cells.Item(9, 17).Value = 21.36443926639834
$ws.Cells.Item(9, 18).Value = 192.279953397585
$ws.Cells.Item(9, 19).Value = 0.02131053200014788
$ws.Cells.Item(9, 20).Value = 0.02131053200014788

